$d = $word.ActiveDocument

# 1. Title text fix: "TTCKT" -> "TCKT" (in the file name mentioned in the heading)
$d.Content.Find.Execute("BI_TLTL_TTCKT_NAM_QUY.xlsx", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "BI_TLTL_TCKT_NAM_QUY.xlsx", 2)

# 2. Merge the runs "B" "2" ":" " Nhap" " " "c" into a single run's text "B2: Nhap c"
#    The surrounding runs/words stay the same; only the literal text needs to read "B2: Nhap c"
#    instead of "B" + "2" + ":" + " Nhap" + " " + "c" (which already renders identically,
#    so no visible text changes - only run-splitting collapses).
$d.Content.Find.Execute("B2: Nhập c", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "B2: Nhập c", 2)

# 3. Split "cung khong the nhap cac du lieu o ngoai bang" so that the leading "C"
#    becomes its own run, followed by "ung khong the nhap cac du lieu o ngoai bang"
$d.Content.Find.Execute("cũng không thể nhập các dữ liệu ở ngoài bảng", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Cũng không thể nhập các dữ liệu ở ngoài bảng", 2)
